$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values for rows where the final delta-stack (dSF)
# differs from the initial value (dS0) after repulling / pushing all data and
# recalculating the mean. Only column F values changed; all other cells are unchanged.
$ws.Range("F13").Value = 0
$ws.Range("F22").Value = 2
$ws.Range("F23").Value = 2
$ws.Range("F26").Value = -2
$ws.Range("F27").Value = 2
$ws.Range("F30").Value = -2
$ws.Range("F35").Value = 0
$ws.Range("F36").Value = 1
$ws.Range("F38").Value = -4
$ws.Range("F39").Value = 2
$ws.Range("F41").Value = -1
$ws.Range("F47").Value = -2
$ws.Range("F51").Value = 4
$ws.Range("F52").Value = -6
$ws.Range("F53").Value = 2
$ws.Range("F54").Value = 5
$ws.Range("F55").Value = 4
$ws.Range("F57").Value = 3
$ws.Range("F58").Value = -2
$ws.Range("F64").Value = 4
$ws.Range("F67").Value = 1
$ws.Range("F68").Value = -3
$ws.Range("F70").Value = -7
$ws.Range("F71").Value = -6
$ws.Range("F72").Value = -5
$ws.Range("F75").Value = 0
$ws.Range("F82").Value = -4
